$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '[Tao%Ai%NULL%0,                 Zhenlu%Yang%NULL%2,                 Zhenlu%Yang%NULL%0,                 Hongyan%Hou%NULL%2,                 Hongyan%Hou%NULL%0,                 Chenao%Zhan%NULL%1,                 Chong%Chen%NULL%1,                 Wenzhi%Lv%NULL%2,                 Wenzhi%Lv%NULL%0,                 Qian%Tao%NULL%2,                 Qian%Tao%NULL%0,                 Ziyong%Sun%NULL%0,                 Ziyong%Sun%NULL%0,                 Liming%Xia%xialiming2017@outlook.com%0]'
$ws.Range("E3").Value = '[Harrison X.%Bai%NULL%0,                 Ben%Hsieh%NULL%1,                 Zeng%Xiong%NULL%2,                 Zeng%Xiong%NULL%0,                 Kasey%Halsey%NULL%2,                 Kasey%Halsey%NULL%0,                 Ji Whae%Choi%NULL%2,                 Ji Whae%Choi%NULL%0,                 Thi My Linh%Tran%NULL%2,                 Thi My Linh%Tran%NULL%0,                 Ian%Pan%NULL%2,                 Ian%Pan%NULL%0,                 Lin-Bo%Shi%NULL%2,                 Lin-Bo%Shi%NULL%0,                 Dong-Cui%Wang%NULL%1,                 Ji%Mei%NULL%2,                 Ji%Mei%NULL%0,                 Xiao-Long%Jiang%NULL%1,                 Qiu-Hua%Zeng%NULL%1,                 Thomas K.%Egglin%NULL%1,                 Ping-Feng%Hu%NULL%2,                 Ping-Feng%Hu%NULL%0,                 Saurabh%Agarwal%NULL%1,                 Fangfang%Xie%NULL%2,                 Fangfang%Xie%NULL%0,                 Sha%Li%NULL%1,                 Terrance%Healey%NULL%2,                 Terrance%Healey%NULL%0,                 Michael K.%Atalay%NULL%1,                 Wei-Hua%Liao%liaoweihua2017@163.com%2,                 Wei-Hua%Liao%liaoweihua2017@163.com%0]'
$ws.Range("E4").Value = '[Zenghui%Cheng%xref no email%0,           Yong%Lu%xref no email%1,           Qiqi%Cao%xref no email%1,           Le%Qin%xref no email%1,           Zilai%Pan%xref no email%1,           Fuhua%Yan%xref no email%1,           Wenjie%Yang%xref no email%1]'
$ws.Range("E5").Value = '[Hyewon%Choi%NULL%0,                 Xiaolong%Qi%NULL%0,                 Xiaolong%Qi%NULL%0,                 Soon Ho%Yoon%yshoka@gmail.com%2,                 Soon Ho%Yoon%yshoka@gmail.com%0,                 Sang Joon%Park%NULL%2,                 Sang Joon%Park%NULL%0,                 Kyung Hee%Lee%NULL%2,                 Kyung Hee%Lee%NULL%0,                 Jin Yong%Kim%NULL%0,                 Jin Yong%Kim%NULL%0,                 Young Kyung%Lee%NULL%2,                 Young Kyung%Lee%NULL%0,                 Hongseok%Ko%NULL%2,                 Hongseok%Ko%NULL%0,                 Ki Hwan%Kim%NULL%2,                 Ki Hwan%Kim%NULL%0,                 Chang Min%Park%NULL%2,                 Chang Min%Park%NULL%0,                 Yun-Hyeon%Kim%NULL%2,                 Yun-Hyeon%Kim%NULL%0,                 Junqiang%Lei%NULL%2,                 Junqiang%Lei%NULL%0,                 Jung Hee%Hong%NULL%2,                 Jung Hee%Hong%NULL%0,                 Hyungjin%Kim%NULL%2,                 Hyungjin%Kim%NULL%0,                 Eui Jin%Hwang%NULL%2,                 Eui Jin%Hwang%NULL%0,                 Seung Jin%Yoo%NULL%2,                 Seung Jin%Yoo%NULL%0,                 Ju Gang%Nam%NULL%2,                 Ju Gang%Nam%NULL%0,                 Chang Hyun%Lee%NULL%2,                 Chang Hyun%Lee%NULL%0,                 Jin Mo%Goo%NULL%2,                 Jin Mo%Goo%NULL%0]'
$ws.Range("E6").Value = '[Xiaofeng%Chen%NULL%0,                 Yanyan%Tang%NULL%1,                 Yongkang%Mo%NULL%1,                 Shengkai%Li%NULL%1,                 Daiying%Lin%NULL%1,                 Zhijian%Yang%NULL%1,                 Zhiqi%Yang%NULL%1,                 Hongfu%Sun%NULL%1,                 Jinming%Qiu%NULL%1,                 Yuting%Liao%NULL%1,                 Jianning%Xiao%NULL%1,                 Xiangguang%Chen%NULL%1,                 Xianheng%Wu%NULL%1,                 Renhua%Wu%NULL%1,                 Zhuozhi%Dai%zhuozhi@ualberta.ca%1]'
$ws.Range("E7").Value = '[Damiano%Caruso%NULL%0,                 Marta%Zerunian%NULL%2,                 Marta%Zerunian%NULL%0,                 Michela%Polici%NULL%2,                 Michela%Polici%NULL%0,                 Francesco%Pucciarelli%NULL%2,                 Francesco%Pucciarelli%NULL%0,                 Tiziano%Polidori%NULL%2,                 Tiziano%Polidori%NULL%0,                 Carlotta%Rucci%NULL%2,                 Carlotta%Rucci%NULL%0,                 Gisella%Guido%NULL%2,                 Gisella%Guido%NULL%0,                 Benedetta%Bracci%NULL%2,                 Benedetta%Bracci%NULL%0,                 Chiara%de Dominicis%NULL%2,                 Chiara%de Dominicis%NULL%0,                 Andrea%Laghi%andrea.laghi@uniroma1.it%2,                 Andrea%Laghi%andrea.laghi@uniroma1.it%0]'
$ws.Range("E8").Value = '[Yuki%Himoto%yukihimoto@gmail.com%0,                 Akihiko%Sakata%NULL%1,                 Mitsuhiro%Kirita%NULL%1,                 Takashi%Hiroi%NULL%1,                 Ken-ichiro%Kobayashi%NULL%1,                 Kenji%Kubo%NULL%1,                 Hyunjin%Kim%NULL%1,                 Azusa%Nishimoto%NULL%1,                 Chikara%Maeda%NULL%1,                 Akira%Kawamura%NULL%1,                 Nobuhiro%Komiya%NULL%1,                 Shigeaki%Umeoka%NULL%1]'
$ws.Range("E9").Value = '[Chunqin%Long%NULL%0,                 Huaxiang%Xu%NULL%1,                 Qinglin%Shen%NULL%1,                 Xianghai%Zhang%NULL%1,                 Bing%Fan%26171381@qq.com%0,                 Chuanhong%Wang%NULL%1,                 Bingliang%Zeng%NULL%1,                 Zicong%Li%NULL%1,                 Xiaofen%Li%NULL%1,                 Honglu%Li%NULL%1]'
$ws.Range("E10").Value = '[Congliang%Miao%NULL%0,                 Mengdi%Jin%NULL%1,                 Li%Miao%NULL%1,                 Xinying%Yang%NULL%1,                 Peng%Huang%NULL%1,                 Huanwen%Xiong%NULL%1,                 Peijie%Huang%NULL%1,                 Qi%Zhao%NULL%1,                 Jiang%Du%NULL%1,                 Jiang%Hong%NULL%1]'
$ws.Range("E11").Value = '[Zeying%Wen%NULL%0,                 Yonge%Chi%NULL%2,                 Yonge%Chi%NULL%0,                 Liang%Zhang%NULL%2,                 Liang%Zhang%NULL%0,                 Huan%Liu%NULL%0,                 Huan%Liu%NULL%0,                 Kun%Du%NULL%1,                 Zhengxing%Li%NULL%2,                 Zhengxing%Li%NULL%0,                 Jie%Chen%NULL%2,                 Jie%Chen%NULL%0,                 Liuhui%Cheng%NULL%2,                 Liuhui%Cheng%NULL%0,                 Daoqing%Wang%wangdaoqing1215@126.com%2,                 Daoqing%Wang%wangdaoqing1215@126.com%0]'
$ws.Range("E12").Value = '[Hui%yang%NULL%1,                 Guoqiang%Sun%NULL%0,                 Fei%Tang%NULL%1,                 Min%Peng%NULL%0,                 Ying%Gao%NULL%1,                 Jing%Peng%NULL%1,                 Hui%Xie%NULL%1,                 Yun%Zhao%NULL%1,                 Zhichun%Jin%NULL%1]'
$ws.Range("E13").Value = '[Wanbo%Zhu%NULL%0,                 Kai%Xie%NULL%0,                 Kai%Xie%NULL%0,                 Hui%Lu%NULL%0,                 Lei%Xu%bayinhexl@126.com%0,                 Shusheng%Zhou%zhouss108@163.com%0,                 Shiyuan%Fang%fangshiyuan2008@126.com%0]'
$ws.Range("E14").Value = '[Anthony%Dangis%NULL%0,                 Christopher%Gieraerts%NULL%2,                 Christopher%Gieraerts%NULL%0,                 Yves%De Bruecker%NULL%2,                 Yves%De Bruecker%NULL%0,                 Lode%Janssen%NULL%2,                 Lode%Janssen%NULL%0,                 Hanne%Valgaeren%NULL%2,                 Hanne%Valgaeren%NULL%0,                 Dagmar%Obbels%NULL%2,                 Dagmar%Obbels%NULL%0,                 Marc%Gillis%NULL%2,                 Marc%Gillis%NULL%0,                 Marc%Van Ranst%NULL%0,                 Marc%Van Ranst%NULL%0,                 Johan%Frans%NULL%0,                 Annick%Demeyere%NULL%2,                 Annick%Demeyere%NULL%0,                 Rolf%Symons%rolf.symons@imelda.be%0,                 Rolf%Symons%rolf.symons@imelda.be%0]'
